# Daily auto-push: append the newest data row to the tracking sheet.
#
# The sheet is a simple append-only log with columns:
#   A = date ("YYYY/MM/DD", stored as plain text - every existing cell in
#       this column is literal text, not a real Excel date)
#   B = weekday label, e.g. "木" (plain text)
#   C = hour/time value (number)
#   D = ranking value (number)
#
# This run adds one new row for 2025/10/09:
#   A = 2025/10/09   B = 木   C = 11   D = 19

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the first empty row right after the existing data block.
$newRow = $ws.UsedRange.Rows.Count + 1

$dateCell = $ws.Range("A" + $newRow)

# A leading apostrophe forces the text "2025/10/09" to be stored as a
# literal string instead of being auto-parsed into a date serial number
# (which is what a bare assignment of a date-shaped string would do).
# This keeps column A consistent with every prior row, which all hold
# plain text dates rather than real Excel dates.
$dateCell.Value = "'2025/10/09"

# The apostrophe entry flips on the cell's "quote prefix" formatting
# flag; putting the style back to the workbook default removes that
# again so the new cell's formatting matches its neighbours above it.
$dateCell.Style = "Normal"

$ws.Range("B" + $newRow).Value = "木"
$ws.Range("C" + $newRow).Value = 11
$ws.Range("D" + $newRow).Value = 19
